$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the last existing data row down to the new row
$ws.Range("A54:B54").Copy()
$ws.Range("A55:B55").PasteSpecial(-4122) # xlPasteFormats

# Add the new row of data (Date: 1/8/2026, Error Count: 8)
$ws.Range("A55").Value = [DateTime]"2026-01-08"
$ws.Range("B55").Value = 8

# Update the selection to match the new last row
$ws.Range("A55:B55").Select()
